# Generate Report for Archive
# Update Status from "Ready for handoff" to "In Translation" for the files
# that are currently mid-translation (02988f15-... and 096a9ef0-...),
# on the Overview sheet as well as each per-locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 3 and 4 correspond to
#   02988f15-fd21-4b31-9053-2c39ef0cea9a.md (row 3)
#   096a9ef0-4141-40f2-a2d0-67a78404e9b4.md (row 4)
# Columns: B = zh-cn status, C = de-de status
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

# zh-cn sheet: column C is "Status" for rows 3 and 4
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# de-de sheet: column C is "Status" for rows 3 and 4
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
